# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# Changes:
# 1. On the "Metadata" sheet, set cell B4 (the "Name" row's value) to
#    "AttributionparticuliereVs" (this adds a new shared string).
# 2. Update the "Date" row's value (B8) to the new generation timestamp
#    "2025-07-18T06:40:38+00:00".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "AttributionparticuliereVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
